$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1070
$endRow = 1169

$dateArr = @(45166,45166,44336,44336,44336,44432,44432,44901,44901,44217,44217,44767,44767,44599,44599,44659,44659,44459,44459,44924,44924,44880,44880,45117,45117,44637,44637,45083,44362,44362,45063,45063,44792,44792,44893,44557,44557,44516,44516,44242,44242,44848,44635,44635,44385,44385,44678,44678,44194,44194,44237,44237,45037,45037,44771,44771,45051,45051,44383,44383,45015,45015,44530,44530,44784,44784,44433,44433,44813,44813,44571,44571,44638,44638,44572,44572,44271,44271,44251,44251,44286,44711,44711,44343,44754,44754,44390,44390,44490,44490,44790,44790,44769,44769,44890,44890,44260,44260,44964,44964)
$qualArr = @('Pintón','Primera Pintón','Maduro','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Pintón','Primera Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón','Pintón','Primera Pintón')
$volArr = @(1000,800,1000,1050,420,500,500,1000,600,800,320,1050,320,1050,500,500,350,700,400,1000,800,600,1500,800,500,500,300,1800,500,260,800,780,820,500,500,800,500,250,600,300,500,1850,950,520,800,500,850,450,300,500,200,400,1000,800,800,500,1050,480,850,400,1000,800,800,500,850,300,850,600,840,480,800,500,1050,540,1000,320,800,360,560,340,800,1050,840,600,400,400,450,300,500,300,1000,850,800,450,500,300,800,400,200,400)
$minArr = @(17000,18000,7000,10000,12000,12000,13000,28000,30000,13000,14000,31000,32000,12000,13000,13000,14000,17000,18000,18000,19000,28000,30000,13000,15000,17000,18000,15000,12000,14000,13000,14000,19000,20000,32000,12000,13000,15000,14000,9000,10000,25000,17000,18000,9000,10000,10000,13000,13000,12000,12000,14000,18000,19000,30000,31000,14000,16000,9000,11000,18000,19000,17000,18000,19000,20000,12000,14000,18000,19000,15000,16000,18500,19500,15000,16000,12000,13000,9000,10000,13000,14000,15000,12000,20000,22000,12000,10000,23000,24000,20000,21000,30000,31000,26000,27000,13000,15000,24000,25000)
$maxArr = @(17000,18000,7000,10000,12000,12000,13000,28000,30000,13000,14000,31000,32000,12000,13000,13000,14000,17000,18000,18000,19000,28000,30000,13000,15000,17000,18000,16000,12000,14000,13000,15000,19000,20000,32000,12000,13000,15000,14000,9000,10000,26000,17000,18000,9000,10000,11000,13000,13000,12000,12000,14000,18000,20000,30000,31000,15000,16000,9000,11000,18000,19000,17000,18000,19000,20000,12000,14000,18000,19000,15000,16000,18500,19500,15000,16000,12000,13000,9000,10000,13000,14000,15000,12000,20000,22000,12000,10000,23000,24000,20000,21000,30000,31000,26000,27000,13000,15000,24000,25000)
$avgArr = @(17000,18000,7000,10000,12000,12000,13000,28000,30000,13000,14000,31000,32000,12000,13000,13000,14000,17000,18000,18000,19000,28000,30000,13000,15000,17000,18000,15444,12000,14000,13000,14641,19000,20000,32000,12000,13000,15000,14000,9000,10000,25459,17000,18000,9000,10000,10765,13000,13000,12000,12000,14000,18000,19375,30000,31000,14524,16000,9000,11000,18000,19000,17000,18000,19000,20000,12000,14000,18000,19000,15000,16000,18500,19500,15000,16000,12000,13000,9000,10000,13000,14000,15000,12000,20000,22000,12000,10000,23000,24000,20000,21000,30000,31000,26000,27000,13000,15000,24000,25000)
$kgArr = @(850,900,350,500,600,600,650,1400,1500,650,700,1550,1600,600,650,650,700,850,900,900,950,1400,1500,650,750,850,900,772,600,700,650,732,950,1000,1600,600,650,750,700,450,500,1273,850,900,450,500,538,650,650,600,600,700,900,969,1500,1550,726,800,450,550,900,950,850,900,950,1000,600,700,900,950,750,800,925,975,750,800,600,650,450,500,650,700,750,600,1000,1100,600,500,1150,1200,1000,1050,1500,1550,1300,1350,650,750,1200,1250)
$originArr = @('Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Bolivia','Bolivia','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador','Ecuador')

for ($i = 0; $i -lt $dateArr.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 4).Value = $dateArr[$i]       # D Fecha
    $ws.Cells.Item($r, 12).Value = $qualArr[$i]      # L Calidad
    $ws.Cells.Item($r, 13).Value = $volArr[$i]       # M Volumen
    $ws.Cells.Item($r, 14).Value = $minArr[$i]       # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $maxArr[$i]       # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $avgArr[$i]       # P Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value = $originArr[$i]    # R Origen
    $ws.Cells.Item($r, 19).Value = $kgArr[$i]        # S Precio $/Kg
}

# Fill constant template columns + number format for the two brand-new rows (1168, 1169)
foreach ($r in @(1168,1169)) {
    $ws.Cells.Item($r, 1).Value = 5                                   # A Mercado ID
    $ws.Cells.Item($r, 2).Value = 'Macroferia Regional de Talca'       # B Mercado
    $ws.Cells.Item($r, 3).Value = 'Maule'                             # C Region
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(1166, 4).NumberFormat  # D style
    $ws.Cells.Item($r, 5).Value = 7                                   # E Codreg
    $ws.Cells.Item($r, 6).Value = 'Fruta'                             # F Tipo
    $ws.Cells.Item($r, 7).Value = 100108                             # G Producto ID
    $ws.Cells.Item($r, 8).Value = 'Tropicales y subtropicales'        # H Producto
    $ws.Cells.Item($r, 9).Value = 100108006                          # I Categoria ID
    $ws.Cells.Item($r, 10).Value = 'Plátano'                          # J Categoria
    $ws.Cells.Item($r, 11).Value = 'Sin especificar'                 # K Variedad
    $ws.Cells.Item($r, 17).Value = '$/caja 20 kilos'                 # Q Unidad de comercializacion
    $ws.Cells.Item($r, 20).Value = 20                                 # T Kg / unidad
}